$d = $word.ActiveDocument

# 1. Merge "Developed and deployed software applications for intelligent " +
#    "vehicle routing" + " optimization and supply chain network design. "
#    into a single run by re-confirming the text via Find & Replace (same
#    text, same formatting) -- Word coalesces adjacent runs with identical
#    formatting when the text is re-set this way.
$d.Content.Find.Execute("vehicle routing", $true, $false, $false, $false, $false, $true, 1, $false, "vehicle routing", 2) | Out-Null

# 2. Change "or" to "and" in "... on premise or in the cloud ..." -- locate
#    the unique phrase, compute the exact sub-range of the word "or", and
#    replace just that word so Word splits the run into three pieces
#    (" on premise ", "and", " in the cloud...") while keeping identical
#    run formatting on all three (toggling a format on/off forces Word to
#    keep the split instead of re-coalescing the runs on save).
$ctx = $d.Content
$found = $ctx.Find.Execute("premise or in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $orStart = $ctx.Start + 8
    $orEnd = $orStart + 2
    $orRange = $d.Range($orStart, $orEnd)
    $orRange.Text = "and"
    $orRange.Font.Bold = 1
    $orRange.Font.Bold = 0
}

# 3. Delete the trailing empty paragraph (pBdr all-nil, rPr color only) that
#    sits right before the final sectPr, leaving the Publications hyperlink
#    paragraph as the last paragraph in the body.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$killRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$killRange.Delete()
